$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("ID", "FORM", "LEMMA", "UPOS", "XPOS", "FEATS", "HEAD", "DEPREL", "DEPS", "MISC", "dist")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 5 + $i  # column E = 5
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}

# Copy the header style (bold, centered, bordered) from D1 to E1:O1
$ws.Range("D1").Copy()
$ws.Range("E1:O1").PasteSpecial(-4122)  # xlPasteFormats
